$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values are written as a leading-apostrophe string so Excel keeps them
# as text (matching the source data, which stores numbers like "597.71"
# or "1.00" as literal text, not numeric values) without altering cell
# formatting/styles.
$ws.Range("D2").Value = "68.052.40"
$ws.Range("E2").Value = "  -0.83%  "
$ws.Range("D3").Value = "3.780.24"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'597.71"
$ws.Range("E5").Value = "  -0.71%  "
$ws.Range("D6").Value = "'170.24"
$ws.Range("E6").Value = "  +0.64%  "
$ws.Range("D7").Value = "3.780.16"
$ws.Range("E7").Value = "  -1.74%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("E10").Value = "  -1.25%  "
$ws.Range("D11").Value = "'6.53"
$ws.Range("E11").Value = "  +0.55%  "
$ws.Range("D12").Value = "'0.454"
$ws.Range("E12").Value = "  -1.34%  "
$ws.Range("E13").Value = "  +1.47%  "
$ws.Range("D14").Value = "'36.64"
$ws.Range("E14").Value = "  -0.63%  "
$ws.Range("D15").Value = "4.416.96"
$ws.Range("E15").Value = "  -1.63%  "
$ws.Range("D16").Value = "3.776.00"
$ws.Range("E16").Value = "  -3.10%  "
$ws.Range("D17").Value = "'18.90"
$ws.Range("E17").Value = "  +3.10%  "
$ws.Range("D18").Value = "68.024.39"
$ws.Range("E18").Value = "  -0.87%  "
$ws.Range("D19").Value = "'7.21"
$ws.Range("E19").Value = "  -1.75%  "
$ws.Range("E20").Value = "  +0.66%  "
$ws.Range("D21").Value = "'10.63"
$ws.Range("E21").Value = "  -2.74%  "
$ws.Range("D22").Value = "'468.60"
$ws.Range("E22").Value = "  -0.85%  "
$ws.Range("D23").Value = "'0.720"
$ws.Range("E23").Value = "  -0.66%  "
$ws.Range("E24").Value = "  -7.47%  "
$ws.Range("D25").Value = "'83.77"
$ws.Range("E25").Value = "  +0.56%  "
$ws.Range("E26").Value = "  +0.30%  "
$ws.Range("D27").Value = "'12.16"
$ws.Range("E27").Value = "  +0.67%  "
$ws.Range("D28").Value = "'10.54"
$ws.Range("E28").Value = "  +1.18%  "
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("E30").Value = "  -1.02%  "
$ws.Range("D31").Value = "3.930.87"
$ws.Range("E31").Value = "  -1.63%  "
$ws.Range("D32").Value = "'7.61"
$ws.Range("E32").Value = "  -1.36%  "
$ws.Range("D33").Value = "'30.53"
$ws.Range("E33").Value = "  -2.60%  "
$ws.Range("D34").Value = "'2.24"
$ws.Range("E34").Value = "  -2.05%  "
$ws.Range("D35").Value = "'9.24"
$ws.Range("E35").Value = "  -0.31%  "
$ws.Range("D36").Value = "3.740.43"
$ws.Range("E36").Value = "  -1.85%  "
$ws.Range("E37").Value = "  -0.91%  "
$ws.Range("E38").Value = "  +0.63%  "
$ws.Range("D39").Value = "'0.139"
$ws.Range("E39").Value = "  -0.27%  "
$ws.Range("D40").Value = "'1.00"
$ws.Range("E40").Value = "  -1.43%  "
$ws.Range("D41").Value = "'5.85"
$ws.Range("E41").Value = "  -1.10%  "
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("E43").Value = "  -0.18%  "
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").Value = "'1.96"
$ws.Range("E45").Value = "  -1.96%  "
$ws.Range("B46").Value = "Cosmos"
$ws.Range("C46").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D46").Value = "'8.69"
$ws.Range("E46").Value = "  +0.39%  "
$ws.Range("D47").Value = "'405.14"
$ws.Range("E47").Value = "  -3.51%  "
$ws.Range("D48").Value = "'0.000284"
$ws.Range("E48").Value = "  -5.32%  "
$ws.Range("D49").Value = "'45.67"
$ws.Range("E49").Value = "  -2.62%  "
$ws.Range("D50").Value = "'40.17"
$ws.Range("E50").Value = "  +7.16%  "
$ws.Range("D51").Value = "'140.92"
$ws.Range("E51").Value = "  -0.56%  "
